# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 3.85
$ws.Range("G2").Value = 4.2
$ws.Range("H2").Value = 2.12
$ws.Range("I2").Value = 2.26
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 3.3
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 1.8
$ws.Range("Q2").Value = 2.06
$ws.Range("R2").Value = 1.31
$ws.Range("S2").Value = 3.7
$ws.Range("T2").Value = 1.83
$ws.Range("U2").Value = 2.04
$ws.Range("V2").Value = 1.79
$ws.Range("W2").Value = 1.31
$ws.Range("X2").Value = 970
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 970
$ws.Range("AA2").Value = 970
$ws.Range("AB2").Value = 16.5
$ws.Range("AC2").Value = 7.8
$ws.Range("AD2").Value = 13
$ws.Range("AE2").Value = 970
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 970
$ws.Range("AH2").Value = 970
$ws.Range("AI2").Value = 970
$ws.Range("AK2").Value = 65
$ws.Range("AL2").Value = 65
$ws.Range("AM2").Value = 130
$ws.Range("AN2").Value = 60
$ws.Range("AO2").Value = 970

# Row 3
$ws.Range("F3").Value = 1.44
$ws.Range("G3").Value = 1.59
$ws.Range("H3").Value = 7.4
$ws.Range("I3").Value = 10.5
$ws.Range("J3").Value = 4.3
$ws.Range("K3").Value = 5.2
$ws.Range("P3").Value = 2.02
$ws.Range("Q3").Value = 1.77

# Row 4
$ws.Range("F4").Value = 1.96
$ws.Range("G4").Value = 1.98
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 3.75
$ws.Range("AC4").Value = 8.4
$ws.Range("AE4").Value = 55
$ws.Range("AO4").Value = 55

# Row 5
$ws.Range("F5").Value = 1.33
$ws.Range("G5").Value = 1.35
$ws.Range("H5").Value = 13
$ws.Range("I5").Value = 14
$ws.Range("J5").Value = 5.4
$ws.Range("T5").Value = 2.46

# Row 6
$ws.Range("F6").Value = 2.7
$ws.Range("G6").Value = 3.4
$ws.Range("H6").Value = 2.48
$ws.Range("I6").Value = 2.82
$ws.Range("J6").Value = 3.35
$ws.Range("K6").Value = 3.95
$ws.Range("P6").Value = 1.98
$ws.Range("Q6").Value = 1.83

# Row 7
$ws.Range("F7").Value = 1.65
$ws.Range("G7").Value = 1.94
$ws.Range("H7").Value = 5
$ws.Range("J7").Value = 3.45
$ws.Range("K7").Value = 4.6
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 3.65
$ws.Range("P7").Value = 2.06
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 1.42
$ws.Range("S7").Value = 2.88
$ws.Range("T7").Value = 1.75
$ws.Range("U7").Value = 2.04
$ws.Range("X7").Value = 22
$ws.Range("Y7").Value = 25
$ws.Range("AB7").Value = 11.5
$ws.Range("AC7").Value = 11.5
$ws.Range("AD7").Value = 26
$ws.Range("AF7").Value = 13.5
$ws.Range("AG7").Value = 12.5
$ws.Range("AH7").Value = 24
$ws.Range("AJ7").Value = 21
$ws.Range("AK7").Value = 21
$ws.Range("AL7").Value = 40
$ws.Range("AN7").Value = 11.5

# Row 8
$ws.Range("F8").Value = 11
$ws.Range("K8").Value = 7.6
$ws.Range("P8").Value = 4.2
$ws.Range("R8").Value = 2.3
$ws.Range("S8").Value = 1.71
$ws.Range("Y8").Value = 21
$ws.Range("Z8").Value = 14
$ws.Range("AC8").Value = 19.5
$ws.Range("AD8").Value = 12.5
$ws.Range("AG8").Value = 42
$ws.Range("AH8").Value = 24
$ws.Range("AI8").Value = 25
$ws.Range("AN8").Value = 80
$ws.Range("AO8").Value = 2.94

# Row 9
$ws.Range("G9").Value = 1.92
$ws.Range("K9").Value = 4.3
$ws.Range("N9").Value = 5.3
$ws.Range("O9").Value = 1.21
$ws.Range("P9").Value = 2.48
$ws.Range("Q9").Value = 1.63
$ws.Range("R9").Value = 1.59
$ws.Range("S9").Value = 2.58
$ws.Range("T9").Value = 1.63
$ws.Range("U9").Value = 2.48
$ws.Range("AI9").Value = 44
$ws.Range("AK9").Value = 17.5
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 8.800000000000001

# Row 10
$ws.Range("N10").Value = 5.3
$ws.Range("S10").Value = 2.64
$ws.Range("T10").Value = 1.66
$ws.Range("U10").Value = 2.44
$ws.Range("Y10").Value = 21
$ws.Range("AD10").Value = 18.5
$ws.Range("AG10").Value = 10.5
$ws.Range("AK10").Value = 18
$ws.Range("AN10").Value = 8.800000000000001

# Row 11
$ws.Range("J11").Value = 9.6
$ws.Range("K11").Value = 10.5
$ws.Range("Q11").Value = 1.55
$ws.Range("S11").Value = 2.42
$ws.Range("U11").Value = 1.53

# Row 12
$ws.Range("F12").Value = 1.44
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 5.1
$ws.Range("P12").Value = 2
$ws.Range("Q12").Value = 1.81

# Row 13
$ws.Range("G13").Value = 7.6
$ws.Range("H13").Value = 1.51
$ws.Range("I13").Value = 1.54
$ws.Range("J13").Value = 4.8
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 1.2
$ws.Range("P13").Value = 2.4
$ws.Range("Q13").Value = 1.6
$ws.Range("R13").Value = 1.55
$ws.Range("S13").Value = 2.54
$ws.Range("T13").Value = 1.76
$ws.Range("U13").Value = 2.1
$ws.Range("X13").Value = 29
$ws.Range("Y13").Value = 13
$ws.Range("Z13").Value = 12.5
$ws.Range("AA13").Value = 17
$ws.Range("AB13").Value = 32
$ws.Range("AC13").Value = 14
$ws.Range("AD13").Value = 12.5
$ws.Range("AE13").Value = 18.5
$ws.Range("AF13").Value = 65
$ws.Range("AG13").Value = 27
$ws.Range("AH13").Value = 26
$ws.Range("AI13").Value = 36
$ws.Range("AK13").Value = 100
$ws.Range("AL13").Value = 90
$ws.Range("AN13").Value = 95
$ws.Range("AO13").Value = 7.4

# Row 14
$ws.Range("F14").Value = 1.36
$ws.Range("G14").Value = 1.44
$ws.Range("H14").Value = 8.800000000000001
$ws.Range("I14").Value = 11
$ws.Range("J14").Value = 4.6
$ws.Range("K14").Value = 5.9
$ws.Range("M14").Value = 1.04
$ws.Range("N14").Value = 4.8
$ws.Range("O14").Value = 1.21
$ws.Range("P14").Value = 2.3
$ws.Range("Q14").Value = 1.65
$ws.Range("R14").Value = 1.44
$ws.Range("S14").Value = 2.36
$ws.Range("T14").Value = 1.94
$ws.Range("U14").Value = 1.88
$ws.Range("V14").Value = 1.1
$ws.Range("W14").Value = 3.25
$ws.Range("X14").Value = 27
$ws.Range("Y14").Value = 40
$ws.Range("Z14").Value = 110
$ws.Range("AA14").Value = 380
$ws.Range("AB14").Value = 11.5
$ws.Range("AC14").Value = 15
$ws.Range("AD14").Value = 42
$ws.Range("AE14").Value = 170
$ws.Range("AF14").Value = 11
$ws.Range("AG14").Value = 12.5
$ws.Range("AH14").Value = 32
$ws.Range("AI14").Value = 140
$ws.Range("AJ14").Value = 14
$ws.Range("AK14").Value = 18
$ws.Range("AL14").Value = 42
$ws.Range("AM14").Value = 170
$ws.Range("AN14").Value = 6.8
